$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6921.5557
$ws.Range("I74").Value = 8165.6665
$ws.Range("J74").Value = 4433.3335
$ws.Range("K74").Value = 8165.6665
$ws.Range("L74").Value = 4433.3335
$ws.Range("M74").Value = -7229.6665
$ws.Range("N74").Value = -6305.3335

$ws.Range("H77").Value = 6921.5557
$ws.Range("I77").Value = 8165.6665
$ws.Range("J77").Value = 4433.3335
$ws.Range("K77").Value = 40828.3325
$ws.Range("L77").Value = 22166.6675
$ws.Range("M77").Value = -36148.3325
$ws.Range("N77").Value = -31526.6675

$ws.Range("H137").Value = 1416.7667
$ws.Range("I137").Value = 1152.7778
$ws.Range("J137").Value = 1812.75
$ws.Range("K137").Value = 3458.3334
$ws.Range("L137").Value = 5438.25
$ws.Range("M137").Value = -908.3334000000004
$ws.Range("N137").Value = -10538.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H61").Value = 1045.091
$ws.Range("I61").Value = 952.8182
$ws.Range("J61").Value = 1321.909
$ws.Range("K61").Value = 952.8182
$ws.Range("L61").Value = 1321.909
$ws.Range("M61").Value = -740.8182
$ws.Range("N61").Value = -1745.909

$ws.Range("H74").Value = 1141
$ws.Range("I74").Value = 1031.375
$ws.Range("J74").Value = 1433.3334
$ws.Range("K74").Value = 1031.375
$ws.Range("L74").Value = 1433.3334
$ws.Range("M74").Value = -157.375
$ws.Range("N74").Value = -3181.3334

$ws.Range("H77").Value = 1141
$ws.Range("I77").Value = 1031.375
$ws.Range("J77").Value = 1433.3334
$ws.Range("K77").Value = 5156.875
$ws.Range("L77").Value = 7166.666999999999
$ws.Range("M77").Value = -788.875
$ws.Range("N77").Value = -15902.667

$ws.Range("H132").Value = 1697.3125
$ws.Range("I132").Value = 1277.4688
$ws.Range("J132").Value = 2537
$ws.Range("K132").Value = 3832.4064
$ws.Range("L132").Value = 7611
$ws.Range("M132").Value = -1302.4064
$ws.Range("N132").Value = -12671

$ws.Range("H136").Value = 1045.091
$ws.Range("I136").Value = 952.8182
$ws.Range("J136").Value = 1321.909
$ws.Range("K136").Value = 2858.4546
$ws.Range("L136").Value = 3965.727
$ws.Range("M136").Value = -308.4546
$ws.Range("N136").Value = -9065.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 16780
$ws.Range("J53").Value = 16780
$ws.Range("L53").Value = 16780
$ws.Range("N53").Value = -17928

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 124.51724
$ws.Range("I7").Value = 127.35
$ws.Range("J7").Value = 118.22222
$ws.Range("K7").Value = 127.35
$ws.Range("L7").Value = 118.22222
$ws.Range("M7").Value = -14.34999999999999
$ws.Range("N7").Value = -344.22222

$ws.Range("H58").Value = 911.6786
$ws.Range("I58").Value = 717.95654
$ws.Range("J58").Value = 1802.8
$ws.Range("K58").Value = 717.95654
$ws.Range("L58").Value = 1802.8
$ws.Range("M58").Value = -514.95654
$ws.Range("N58").Value = -2208.8

$ws.Range("H64").Value = 26000
$ws.Range("J64").Value = 26000
$ws.Range("L64").Value = 26000
$ws.Range("N64").Value = -26496

$ws.Range("H67").Value = 26000
$ws.Range("J67").Value = 26000
$ws.Range("L67").Value = 26000
$ws.Range("N67").Value = -27716

$ws.Range("H122").Value = 1531.0834
$ws.Range("I122").Value = 1376.5
$ws.Range("J122").Value = 1685.6666
$ws.Range("K122").Value = 4129.5
$ws.Range("L122").Value = 5056.9998
$ws.Range("M122").Value = -1679.5
$ws.Range("N122").Value = -9956.9998

$ws.Range("H134").Value = 1109.3405
$ws.Range("I134").Value = 975.8889
$ws.Range("J134").Value = 1546.091
$ws.Range("K134").Value = 2927.6667
$ws.Range("L134").Value = 4638.272999999999
$ws.Range("M134").Value = -392.6667000000002
$ws.Range("N134").Value = -9708.272999999999

$ws.Range("H136").Value = 911.6786
$ws.Range("I136").Value = 717.95654
$ws.Range("J136").Value = 1802.8
$ws.Range("K136").Value = 2153.86962
$ws.Range("L136").Value = 5408.4
$ws.Range("M136").Value = 396.1303800000001
$ws.Range("N136").Value = -10508.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 393.45456
$ws.Range("I15").Value = 84.75
$ws.Range("J15").Value = 569.8570999999999
$ws.Range("K15").Value = 254.25
$ws.Range("L15").Value = 1709.5713
$ws.Range("M15").Value = -114.25
$ws.Range("N15").Value = -1989.5713

$ws.Range("H70").Value = 1583.3334
$ws.Range("I70").Value = 1583.3334
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4750.0002
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4435.0002
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 1583.3334
$ws.Range("I73").Value = 1583.3334
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4750.0002
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3658.0002
$ws.Range("N73").ClearContents()

$ws.Range("H131").Value = 872.7568
$ws.Range("J131").Value = 991.6070999999999
$ws.Range("L131").Value = 2974.8213
$ws.Range("N131").Value = -13054.8213

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 3834.6667
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 3834.6667
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 3834.6667
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -4172.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3233.111
$ws.Range("I122").Value = 3073.1
$ws.Range("J122").Value = 4033.1667
$ws.Range("K122").Value = 9219.299999999999
$ws.Range("L122").Value = 12099.5001
$ws.Range("M122").Value = -6769.299999999999
$ws.Range("N122").Value = -16999.5001

$ws.Range("H132").Value = 1629.409
$ws.Range("I132").Value = 1153.8438
$ws.Range("J132").Value = 2897.5833
$ws.Range("K132").Value = 3461.5314
$ws.Range("L132").Value = 8692.749899999999
$ws.Range("M132").Value = -931.5314000000003
$ws.Range("N132").Value = -13752.7499

$ws.Range("H136").Value = 2113.5293
$ws.Range("I136").Value = 1109.5
$ws.Range("J136").Value = 3547.8572
$ws.Range("K136").Value = 3328.5
$ws.Range("L136").Value = 10643.5716
$ws.Range("M136").Value = -778.5
$ws.Range("N136").Value = -15743.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H122").Value = 1736.3334
$ws.Range("I122").Value = 1537.75
$ws.Range("J122").Value = 2001.1111
$ws.Range("K122").Value = 4613.25
$ws.Range("L122").Value = 6003.3333
$ws.Range("M122").Value = -2163.25
$ws.Range("N122").Value = -10903.3333
